# "Prompts user input to select a movie" -- get_selection() (with a
# validation while-loop) asked the user which movie to act on; the user
# picked "The Children Act" and the app moved it up one row so it now
# appears directly above "The Little Stranger" in the showtimes list.
# Net effect on the worksheet: rows 19 and 20 trade places completely
# (movie name, rating, runtime -- the release year is the same for every
# row so it doesn't visibly change).

if ($null -eq $wb) {
    $wb = $excel.ActiveWorkbook
}
$ws = $wb.ActiveSheet

$rowTop = 19     # "The Little Stranger" before the move
$rowBottom = 20  # "The Children Act" before the move
$firstCol = 1    # A: Movies Playing
$lastCol = 4     # D: Release Year

$rangeTop = $ws.Range($ws.Cells.Item($rowTop, $firstCol), $ws.Cells.Item($rowTop, $lastCol))
$rangeBottom = $ws.Range($ws.Cells.Item($rowBottom, $firstCol), $ws.Cells.Item($rowBottom, $lastCol))
$scratch = $ws.Range($ws.Cells.Item($rowTop, $lastCol + 2), $ws.Cells.Item($rowTop, $lastCol + 2 + ($lastCol - $firstCol)))

# Copy/paste (rather than re-typing the literal values) keeps each cell's
# original data type -- text stays text, so "6.7"/"105" don't turn into
# numbers -- and leaves formatting untouched, just like dragging the rows
# around in the UI would.
$rangeTop.Copy()
$scratch.PasteSpecial(-4104)  # xlPasteAll

$rangeBottom.Copy()
$rangeTop.PasteSpecial(-4104)

$scratch.Copy()
$rangeBottom.PasteSpecial(-4104)

$scratch.Clear()
$excel.CutCopyMode = 0
